$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Summary sheet ("总计"): insert a new row 2 for 2022-Q3, pushing
#    the existing quarterly rows down by one.
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Shift A2:D2 down (and everything below) to make room for the new row.
$summary.Range("A2:D2").Insert(-4121)  # xlShiftDown

# Copy the formatting of the row that is now pushed down (row 3, the
# old "2022-Q2" row) onto the freshly inserted row 2 so the new row
# matches the existing look (bold/bordered index cell in column A).
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$summary.Range("B3:D3").Copy()
$summary.Range("B2:D2").PasteSpecial(-4122)  # xlPasteFormats

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.16

# ------------------------------------------------------------------
# 2) New quarterly sheet "2022-Q3": clone the "2022-Q2" sheet (so it
#    inherits the same column layout / styles) and place it right
#    before "2022-Q2", then overwrite its data with the Q3 figures.
# ------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q2Sheet.Copy($q2Sheet, $null)
$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# The template ("2022-Q2") has 4 data rows; Q3 only needs 2, so drop
# the extra two rows (and their trailing blank dimension).
$q3Sheet.Range("A4:H5").Delete()

function Set-TextCell($rng, [string]$val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell $q3Sheet.Range("B2") "515760"
Set-TextCell $q3Sheet.Range("C2") "华夏中证浙江国资创新发展ETF"
Set-TextCell $q3Sheet.Range("D2") "2.04"
Set-TextCell $q3Sheet.Range("E2") "99.57"
Set-TextCell $q3Sheet.Range("F2") "6.23"
Set-TextCell $q3Sheet.Range("G2") "0.1271"
$q3Sheet.Range("H2").Value = 4

Set-TextCell $q3Sheet.Range("B3") "512190"
Set-TextCell $q3Sheet.Range("C3") "浙商汇金中证浙江凤凰行动50ETF"
Set-TextCell $q3Sheet.Range("D3") "0.48"
Set-TextCell $q3Sheet.Range("E3") "98.92"
Set-TextCell $q3Sheet.Range("F3") "6.61"
Set-TextCell $q3Sheet.Range("G3") "0.0317"
$q3Sheet.Range("H3").Value = 3
